$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 690 ("Vega Modelo de Temuco - Ají"
# weekly price table). This pushes the existing rows 690-700 down to 693-703
# and leaves three blank rows (690-692) ready to receive this week's data.
$ws.Rows("690:692").Insert()

# Row 690: new weekly record
$ws.Range("A690").Value = 10
$ws.Range("B690").Value = "Vega Modelo de Temuco"
$ws.Range("C690").Value = "La Araucanía"
$ws.Range("D690").Value = 44656
$ws.Range("E690").Value = 9
$ws.Range("F690").Value = 100112021
$ws.Range("G690").Value = "Ají"
$ws.Range("H690").Value = "Americana (o)"
$ws.Range("I690").Value = "Primera"
$ws.Range("J690").Value = 65
$ws.Range("K690").Value = 17000
$ws.Range("L690").Value = 17000
$ws.Range("M690").Value = 17000
$ws.Range("N690").Value = "`$/caja 15 kilos"
$ws.Range("O690").Value = "Región del Maule"
$ws.Range("P690").Value = 1133
$ws.Range("Q690").Value = 15
$ws.Range("R690").Value = "Hortaliza"

# Row 691: new weekly record
$ws.Range("A691").Value = 10
$ws.Range("B691").Value = "Vega Modelo de Temuco"
$ws.Range("C691").Value = "La Araucanía"
$ws.Range("D691").Value = 44656
$ws.Range("E691").Value = 9
$ws.Range("F691").Value = 100112021
$ws.Range("G691").Value = "Ají"
$ws.Range("H691").Value = "Chilena(o)"
$ws.Range("I691").Value = "Primera"
$ws.Range("J691").Value = 110
$ws.Range("K691").Value = 15000
$ws.Range("L691").Value = 15000
$ws.Range("M691").Value = 15000
$ws.Range("N691").Value = "`$/caja 15 kilos"
$ws.Range("O691").Value = "Región del Maule"
$ws.Range("P691").Value = 1000
$ws.Range("Q691").Value = 15
$ws.Range("R691").Value = "Hortaliza"

# Row 692: new weekly record
$ws.Range("A692").Value = 10
$ws.Range("B692").Value = "Vega Modelo de Temuco"
$ws.Range("C692").Value = "La Araucanía"
$ws.Range("D692").Value = 44656
$ws.Range("E692").Value = 9
$ws.Range("F692").Value = 100112021
$ws.Range("G692").Value = "Ají"
$ws.Range("H692").Value = "Inferno"
$ws.Range("I692").Value = "Extra"
$ws.Range("J692").Value = 110
$ws.Range("K692").Value = 25000
$ws.Range("L692").Value = 25000
$ws.Range("M692").Value = 25000
$ws.Range("N692").Value = "`$/caja 15 kilos"
$ws.Range("O692").Value = "Región de Arica y Parinacota"
$ws.Range("P692").Value = 1667
$ws.Range("Q692").Value = 15
$ws.Range("R692").Value = "Hortaliza"
